$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action1")
$ws.Activate()

$ws.Range("B3").Value = "f2"
$ws.Range("B2").Value = "d2"

$ws.Range("B2").Select()
